# Add a new worksheet "AI Roadmap" at the end of the workbook.
# Copy the last existing sheet first so the new sheet inherits the same
# sheet-level properties (sheetPr/outlinePr, pageMargins, etc.) used
# consistently across the rest of this workbook, then wipe its contents.
$wb = $excel.ActiveWorkbook
$lastIndex = $wb.Worksheets.Count
$srcSheet = $wb.Worksheets.Item($lastIndex)
$srcSheet.Copy($null, $srcSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "AI Roadmap"
$ws.Cells.Clear()

$ws.Range("A1").Value = 'AI ENHANCEMENT ROADMAP (Grok-verified Jan 2026)'
$ws.Range("A3").Value = 'TIER 1: QUICK WINS (1-2 weeks)'
$ws.Range("A4").Value = 'Feature'
$ws.Range("B4").Value = 'Technology'
$ws.Range("C4").Value = 'Impact'
$ws.Range("D4").Value = 'Cost/1K users'
$ws.Range("A5").Value = 'LLM Resume Parsing'
$ws.Range("B5").Value = 'GPT-4/Claude'
$ws.Range("C5").Value = 'Better profile extraction'
$ws.Range("D5").Value = '$2-5'
$ws.Range("A6").Value = 'LLM Job Analysis'
$ws.Range("B6").Value = 'GPT-4/Claude'
$ws.Range("C6").Value = 'True requirements extraction'
$ws.Range("D6").Value = '$2-5'
$ws.Range("A7").Value = 'Natural Language Explanations'
$ws.Range("B7").Value = 'GPT-4o-mini'
$ws.Range("C7").Value = 'Plain English match reasons'
$ws.Range("D7").Value = '$5-10'
$ws.Range("A9").Value = 'TIER 2: COMPETITIVE ADVANTAGES (1-2 months)'
$ws.Range("A10").Value = 'Feature'
$ws.Range("B10").Value = 'Technology'
$ws.Range("C10").Value = 'Impact'
$ws.Range("D10").Value = 'Cost/1K users'
$ws.Range("A11").Value = 'Career Path Prediction'
$ws.Range("B11").Value = 'LLM + labor data'
$ws.Range("C11").Value = 'Strategic career guidance'
$ws.Range("D11").Value = '$5-10'
$ws.Range("A12").Value = 'Skill Gap Analysis'
$ws.Range("B12").Value = 'LLM analysis'
$ws.Range("C12").Value = 'Actionable improvement steps'
$ws.Range("D12").Value = '$3-5'
$ws.Range("A13").Value = 'Personalized Weights'
$ws.Range("B13").Value = 'Reinforcement learning'
$ws.Range("C13").Value = 'Self-improving algorithm'
$ws.Range("D13").Value = '$0 (compute)'
$ws.Range("A14").Value = 'LLM Coach Assistant'
$ws.Range("B14").Value = 'Claude Haiku'
$ws.Range("C14").Value = 'Scale coaching 10x'
$ws.Range("D14").Value = '$10-20'
$ws.Range("A16").Value = 'TIER 3: MOONSHOTS (3-6 months)'
$ws.Range("A17").Value = 'Feature'
$ws.Range("B17").Value = 'Technology'
$ws.Range("C17").Value = 'Impact'
$ws.Range("D17").Value = 'Cost/1K users'
$ws.Range("A18").Value = 'Interview Simulation'
$ws.Range("B18").Value = 'LLM + speech'
$ws.Range("C18").Value = 'Premium coaching feature'
$ws.Range("D18").Value = '$20-50'
$ws.Range("A19").Value = 'Bias Detection'
$ws.Range("B19").Value = 'AI Fairness 360'
$ws.Range("C19").Value = 'Compliance, trust'
$ws.Range("D19").Value = '$5-10'
$ws.Range("A20").Value = 'Market Intelligence'
$ws.Range("B20").Value = 'Real-time APIs'
$ws.Range("C20").Value = 'Strategic positioning'
$ws.Range("D20").Value = '$10-20'
$ws.Range("A22").Value = 'TOTAL AI COST PROJECTION'
$ws.Range("A23").Value = 'Tier 1 (MVP+)'
$ws.Range("B23").Value = '$9-20/1K users/month'
$ws.Range("A24").Value = 'Tier 1+2 (Growth)'
$ws.Range("B24").Value = '$17-35/1K users/month'
$ws.Range("A25").Value = 'Full Stack'
$ws.Range("B25").Value = '$50-100/1K users/month'
$ws.Range("A27").Value = 'ROI AT $19/USER AVG'
$ws.Range("A28").Value = 'Revenue per 1K users'
$ws.Range("B28").Value = '$19,000/month'
$ws.Range("A29").Value = 'AI cost (Tier 1+2)'
$ws.Range("B29").Value = '$35/month'
$ws.Range("A30").Value = 'AI as % of revenue'
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = '0.18%'
